$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "26.999.98"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E2").Value = "  -0.15%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.874.45"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.55%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "305.53"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9998"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +0.06%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.5066"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.3666"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -2.19%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.07208"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +0.77%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.8947"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.88%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "20.74"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.870.44"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +0.26%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07524"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -0.71%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "94.97"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.12%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "5.246"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E15").Value = "  -1.21%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "1.000"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.000008530"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.12%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "14.25"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +1.07%  "
$ws.Range("E19").Value = "  +0.08%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "27.024.29"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.24%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.026"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -0.17%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.093.03"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.47%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "10.40"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -1.22%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.404"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -0.96%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "148.36"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.39%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.786"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -3.05%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "17.89"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -0.43%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.082"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -1.22%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "113.25"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.44%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "4.713"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.83%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "4.699"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.12%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.09159"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.69%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.05122"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -0.38%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7518"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E34").Value = "  +3.01%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "2.983"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -2.25%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "1.161"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +0.54%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.222"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +6.18%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.567"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.86%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.5685"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +7.06%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.02004"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E40").Value = "  -1.78%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "1.074"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -0.03%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "6.627"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +0.90%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "115.62"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -0.71%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.536"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +3.10%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.1476"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.42%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.4758"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +2.87%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.9995"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +0.07%  "
$ws.Range("E48").Value = "  +0.88%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "1.568"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -0.03%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "36.86"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "63.16"
$c.NumberFormat = "General"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -1.20%  "
